# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the f21e88cf handoff row
# on both the zh-cn and de-de worksheets to reflect the newly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-07 08:39:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-07 08:39:23"
